# 2021 SAT Scores.xlsx - apply "Add files via upload" revision
#
# Semantic summary of the change (derived from the OOXML diff):
#   1. A new title row is inserted above the existing header row. The new
#      A1:F1 range is merged and shows the number 2021, centered, with a
#      thin box border. All existing rows shift down by one (old row 1
#      becomes row 2, old row 19 becomes row 20).
#   2. The two "Seniors who had taken the SATs in 2021, ..." column
#      headers drop the "in 2021" wording (now redundant because of the
#      new title row) -> "Seniors who had taken the SATs, ...".
#   3. The print areas / named ranges grow/shift to match the new used
#      range (A1:U19 -> A2:U20, A1:F20 -> A2:F21).
#   4. The active selection becomes C12 (in the new row numbering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new title row at the very top -----------------------
$ws.Rows.Item(1).Insert()

$title = $ws.Range("A1:F1")
$title.Merge()
$title.Value = 2021
$title.HorizontalAlignment = -4108   # xlCenter
$title.Borders.LineStyle = 1         # xlContinuous (thin box around range)

# --- 2. Drop "in 2021" from the two SAT-header labels -----------------
$nl = [char]10
$ws.Range("B2").Value = "Seniors who had taken the SATs, Number" + $nl + "(in thousands) "
$ws.Range("C2").Value = "Seniors who had taken the SATs, Percentage distribution"

# --- 3. Update print areas / named ranges to the new used range -------
$ws.PageSetup.PrintArea = "`$A`$2:`$U`$20"

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Print_Area") {
        $n.RefersTo = "='2021 SAT Scores'!`$A`$2:`$F`$21"
    }
    if ($n.Name -eq "2021 SAT Scores!Print_Area_MI") {
        $n.RefersTo = "='2021 SAT Scores'!`$A`$2:`$F`$21"
    }
}

# --- 4. Restore the active selection on the new layout -----------------
$ws.Range("C12").Select()
